$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The distinctive number-format that currently lives on D10 ("Iguazu
# Falls" row) needs to travel with that row's data down to D13 once the
# table below gets re-ordered - move the formatting first, before the
# values change, swapping D10 back to the plain column-D look.
$ws.Range("D10").Copy()
$ws.Range("D13").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("D10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# The two cities whose population figures were revised upward...
# ...combined with a re-sort of the whole table (rows 2-13) by the
# Population column, descending, produces this final row order.
# (City, Country, Population, Area)

$ws.Range("A2").Value = "Buenos Aires"
$ws.Range("B2").Value = "Argentina"
$ws.Range("C2").Value = 2891000.0
$ws.Range("D2").Value = 4758.0

$ws.Range("A3").Value = "Toronto"
$ws.Range("B3").Value = "Canada"
$ws.Range("C3").Value = 2800000.0
$ws.Range("D3").Value = 2731571.0

$ws.Range("A4").Value = "Pyeongchang"
$ws.Range("B4").Value = "South Korea"
$ws.Range("C4").Value = 2581000.0
$ws.Range("D4").Value = 3194.0

$ws.Range("A5").Value = "Marakesh"
$ws.Range("B5").Value = "Morocco"
$ws.Range("C5").Value = 928850.0
$ws.Range("D5").Value = 200.0

$ws.Range("A6").Value = "Albuquerque"
$ws.Range("B6").Value = "New Mexico"
$ws.Range("C6").Value = 559277.0
$ws.Range("D6").Value = 491.0

$ws.Range("A7").Value = "Los Cabos"
$ws.Range("B7").Value = "Mexico"
$ws.Range("C7").Value = 287651.0
$ws.Range("D7").Value = 3750.0

$ws.Range("A8").Value = "Greenville"
$ws.Range("B8").Value = "USA"
$ws.Range("C8").Value = 84554.0
$ws.Range("D8").Value = 68.0

$ws.Range("A9").Value = "Archipelago Sea"
$ws.Range("B9").Value = "Finland"
$ws.Range("C9").Value = 60000.0
$ws.Range("D9").Value = 8300.0

$ws.Range("A10").Value = "Walla Walla Valley"
$ws.Range("B10").Value = "USA"
$ws.Range("C10").Value = 32237.0
$ws.Range("D10").Value = 33.0

$ws.Range("A11").Value = "Salina Island"
$ws.Range("B11").Value = "Italy"
$ws.Range("C11").Value = 4000.0
$ws.Range("D11").Value = 27.0

$ws.Range("A12").Value = "Solta"
$ws.Range("B12").Value = "Croatia"
$ws.Range("C12").Value = 1700.0
$ws.Range("D12").Value = 59.0

$ws.Range("A13").Value = "Iguazu Falls"
$ws.Range("B13").Value = "Argentina"
$ws.Range("C13").Value = 0.0
$ws.Range("D13").Value = 672.0

# Turn on AutoFilter for the table range.
$null = $ws.Range("A1:D13").AutoFilter()

# Excel backs every AutoFilter range with a hidden workbook-scoped
# defined name - register it explicitly so the saved file matches.
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$D`$13")
$filterName.Visible = $false
